$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.993.13"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.868.63"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.67"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5081"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3916"
$ws.Range("E8").Value = "  -2.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08140"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.02"
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.086"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("E12").Value = "  +5.74%  "
$ws.Range("D13").Value = "1.860.85"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.240"
$ws.Range("E14").Value = "  -1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.43"
$ws.Range("E17").Value = "  -5.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001075"
$ws.Range("E18").Value = "  -3.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06364"
$ws.Range("E19").Value = "  -5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.82"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "29.981.57"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.770"
$ws.Range("E23").Value = "  -4.95%  "
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.198"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "2.079.88"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.04"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.84"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.210"
$ws.Range("E29").Value = "  -10.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.53"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.045"
$ws.Range("E31").Value = "  -3.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1030"
$ws.Range("E32").Value = "  -2.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.860"
$ws.Range("E33").Value = "  -3.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.724"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02404"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.188"
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06311"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2131"
$ws.Range("E38").Value = "  -4.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.171"
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.468"
$ws.Range("E40").Value = "  -6.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6268"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.207"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5870"
$ws.Range("E45").Value = "  -4.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.87"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.619"
$ws.Range("E47").Value = "  -3.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.978"
$ws.Range("E48").Value = "  -3.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.29"
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.199"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.145"
$ws.Range("E51").Value = "  -0.98%  "
